# Reorder the "readme" sheet's summary table columns from
#   index, sheet_name, Date, JobNo, Author
# to
#   index, Author, JobNo, sheet_name, Date
# and refresh the Author/JobNo/Date values + the Project Information
# "Date of Analysis" timestamp, per the removal of datamine_functions.py
# (functionality moved into utils.py, re-run produced a new job id/date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")

$lastRow = $ws.Range("A1").End(-4121).Row  # xlDown

# New header order for columns B:E
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "JobNo"
$ws.Range("D1").Value = "sheet_name"
$ws.Range("E1").Value = "Date"

# The new Date column holds a purely-numeric-looking string ("20220422").
# Pre-format it as Text so Excel stores it as a string instead of silently
# re-typing it as a number (matches how the source data was produced).
$ws.Range("E2:E" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $sheetName = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 2).Value = "jovyan"
    $ws.Cells.Item($r, 3).Value = "/c/e"
    $ws.Cells.Item($r, 4).Value = $sheetName
    $ws.Cells.Item($r, 5).Value = "20220422"
}

# Keep the table's ListColumns in sync with the new header labels.
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item(2).Name = "Author"
$table.ListColumns.Item(3).Name = "JobNo"
$table.ListColumns.Item(4).Name = "sheet_name"
$table.ListColumns.Item(5).Name = "Date"

# Update the recorded analysis timestamp on the Project Information sheet.
$piWs = $wb.Worksheets.Item("Project Information")
$piWs.Range("B12").Value = "2022-04-22 12:42:30.496561"
